$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Sheet1 -> List1)
$ws.Name = "List1"

# Wipe the old formula-driven data/styles so we can lay down the new,
# plain literal values (no shared formulas, no cell styles, 9 rows instead
# of 10 - origin/sprite-center rework for the flying objects).
$ws.Range("A1:B10").Clear()

$values = @(
    @(1, 3),
    @(3, 4),
    @(7, 2),
    @(13, 2),
    @(17, 2),
    @(22, 2),
    @(28, 2),
    @(32, 2),
    @(43, 2)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
}

# The "Zarez" (Comma) cell style is no longer used by any cell - drop it
# from the workbook's style collection.
$wb.Styles.Item("Zarez").Delete()

# Move the active selection to A2, matching the saved view state.
$ws.Range("A2").Select()
